$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that currently sits in the
#    "Objectives" heading paragraph at the top of the document.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. In the "Book Catalog Web App" section, the first sentence changes
#    from "...Add a ratings field..." to "...Add a publisher field...".
#    Replace the word "ratings" with "publisher" using Find/Replace so
#    the Selection range lands exactly on the new word.
$sel = $word.Selection
$found = $sel.Find.Execute("ratings", $true, $false, $false, $false, $false, $true, 1, $false, "publisher", 2)

if ($found) {
    # The replaced word is now selected (Selection.Range covers "publisher").
    $newWordRange = $sel.Range

    # 3. Word leaves its (hidden) "_GoBack" bookmark at the location of the
    #    most recent edit, i.e. immediately after the text that was typed.
    #    Re-create that bookmark collapsed right after "publisher". Adding
    #    it across the whole word first guarantees the surrounding run is
    #    split cleanly on both sides (matching how Word itself breaks runs
    #    around edited text), then we move the bookmark to a collapsed
    #    range at the end of the word.
    $d.Bookmarks.Add("_GoBack", $newWordRange)
    $collapsed = $d.Range($newWordRange.End, $newWordRange.End)
    $d.Bookmarks.Add("_GoBack", $collapsed)
}
